# Scheduled runner update: refresh market-board price snapshots (currentAveragePrice*,
# LevePrice*, LeveProfit*) for the affected Leve rows across the crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 1800
$ws.Range("I40").Value = 1750
$ws.Range("J40").Value = 1866.6666
$ws.Range("K40").Value = 1750
$ws.Range("L40").Value = 1866.6666
$ws.Range("M40").Value = -1575
$ws.Range("N40").Value = -2216.6666

# Row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 4527999
$ws.Range("I64").Value = 8549887
$ws.Range("J64").Value = 3375
$ws.Range("K64").Value = 8549887
$ws.Range("L64").Value = 3375
$ws.Range("M64").Value = -8549639
$ws.Range("N64").Value = -3871

# Row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 4527999
$ws.Range("I67").Value = 8549887
$ws.Range("J67").Value = 3375
$ws.Range("K67").Value = 8549887
$ws.Range("L67").Value = 3375
$ws.Range("M67").Value = -8549029
$ws.Range("N67").Value = -5091

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 3709.476
$ws.Range("I74").Value = 3549.9
$ws.Range("J74").Value = 3854.5454
$ws.Range("K74").Value = 3549.9
$ws.Range("L74").Value = 3854.5454
$ws.Range("M74").Value = -2613.9
$ws.Range("N74").Value = -5726.5454

# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 5424.5625
$ws.Range("I76").Value = 3399
$ws.Range("J76").Value = 7000
$ws.Range("K76").Value = 3399
$ws.Range("L76").Value = 7000
$ws.Range("M76").Value = -3084
$ws.Range("N76").Value = -7630

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 3709.476
$ws.Range("I77").Value = 3549.9
$ws.Range("J77").Value = 3854.5454
$ws.Range("K77").Value = 17749.5
$ws.Range("L77").Value = 19272.727
$ws.Range("M77").Value = -13069.5
$ws.Range("N77").Value = -28632.727

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 5424.5625
$ws.Range("I79").Value = 3399
$ws.Range("J79").Value = 7000
$ws.Range("K79").Value = 3399
$ws.Range("L79").Value = 7000
$ws.Range("M79").Value = -2307
$ws.Range("N79").Value = -9184

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1264.907
$ws.Range("I138").Value = 1180.2683
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 3540.8049
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = 1599.1951
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
# Row 3: Skillet Labor / Bronze Skillet
$ws.Range("H3").Value = 7951
$ws.Range("I3").Value = 902.5
$ws.Range("J3").Value = 14999.5
$ws.Range("K3").Value = 902.5
$ws.Range("L3").Value = 14999.5
$ws.Range("M3").Value = -787.5
$ws.Range("N3").Value = -15229.5

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 24002.4
$ws.Range("I122").Value = 37337.332
$ws.Range("K122").Value = 112011.996
$ws.Range("M122").Value = -109561.996

# Row 128: Heading toward Bankruptcy / Manganese Helm of the Falling Dragon
$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960

$ws = $wb.Worksheets.Item("BSM")
# Row 76: Keep Up with the Mechanics / Titanium-barreled Arquebus
$ws.Range("H76").Value = 37632
$ws.Range("J76").Value = 37632
$ws.Range("L76").Value = 37632
$ws.Range("N76").Value = -38262

# Row 79: Unconventional Weaponry (L) / Titanium-barreled Arquebus
$ws.Range("H79").Value = 37632
$ws.Range("J79").Value = 37632
$ws.Range("L79").Value = 37632
$ws.Range("N79").Value = -39816

$ws = $wb.Worksheets.Item("CRP")
# Row 2: In with the New / Bone Harpoon
$ws.Range("H2").Value = 131375
$ws.Range("I2").Value = 250250
$ws.Range("J2").Value = 12500
$ws.Range("K2").Value = 250250
$ws.Range("L2").Value = 12500
$ws.Range("M2").Value = -250137
$ws.Range("N2").Value = -12726

# Row 60: Bowing to Greater Power / Yew Longbow
$ws.Range("H60").Value = 8716
$ws.Range("J60").Value = 9202.4
$ws.Range("L60").Value = 9202.4
$ws.Range("N60").Value = -10224.4

# Row 94: Beech, Please / Beech Lumber
$ws.Range("H94").Value = 3183.2856
$ws.Range("I94").Value = 2326.2
$ws.Range("J94").Value = 3659.4443
$ws.Range("K94").Value = 2326.2
$ws.Range("L94").Value = 3659.4443
$ws.Range("M94").Value = -1875.2
$ws.Range("N94").Value = -4561.4443

$ws = $wb.Worksheets.Item("CUL")
# Row 86: Let's Not Get Sappy / Birch Syrup
$ws.Range("H86").Value = 1270.2
$ws.Range("I86").Value = 425.5
$ws.Range("J86").Value = 1833.3334
$ws.Range("K86").Value = 1276.5
$ws.Range("L86").Value = 5500.0002
$ws.Range("M86").Value = -90.5
$ws.Range("N86").Value = -7872.0002

# Row 89: Luxury Spillover (L) / Birch Syrup
$ws.Range("H89").Value = 1270.2
$ws.Range("I89").Value = 425.5
$ws.Range("J89").Value = 1833.3334
$ws.Range("K89").Value = 3829.5
$ws.Range("L89").Value = 16500.0006
$ws.Range("M89").Value = 2098.5
$ws.Range("N89").Value = -28356.0006

# Row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 7292.1665
$ws.Range("I137").Value = 5986
$ws.Range("J137").Value = 7794.5386
$ws.Range("K137").Value = 17958
$ws.Range("L137").Value = 23383.6158
$ws.Range("M137").Value = -12858
$ws.Range("N137").Value = -33583.6158

$ws = $wb.Worksheets.Item("GSM")
# Row 5: Hora at Me / Bone Hora
$ws.Range("H5").Value = 300
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 13671799
$ws.Range("J80").Value = 2555295.2
$ws.Range("L80").Value = 2555295.2
$ws.Range("N80").Value = -2557291.2

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 13671799
$ws.Range("J83").Value = 2555295.2
$ws.Range("L83").Value = 12776476
$ws.Range("N83").Value = -12786460

# Row 123: Workplace Workout / Ametrine Ring of Fending
$ws.Range("H123").Value = 20163
$ws.Range("J123").Value = 20163
$ws.Range("L123").Value = 20163
$ws.Range("N123").Value = -25063

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 6523.48
$ws.Range("I7").Value = 9511.714
$ws.Range("J7").Value = 5361.3887
$ws.Range("K7").Value = 9511.714
$ws.Range("L7").Value = 5361.3887
$ws.Range("M7").Value = -9399.714
$ws.Range("N7").Value = -5585.3887

# Row 35: No Risk, No Reward / Toadskin Cesti
$ws.Range("H35").Value = 1110
$ws.Range("I35").Value = 1110
$ws.Range("K35").Value = 1110
$ws.Range("M35").Value = -774

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 1293.619
$ws.Range("I61").Value = 1204.5294
$ws.Range("J61").Value = 1672.25
$ws.Range("K61").Value = 1204.5294
$ws.Range("L61").Value = 1672.25
$ws.Range("M61").Value = -1002.5294
$ws.Range("N61").Value = -2076.25

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1293.619
$ws.Range("I113").Value = 1204.5294
$ws.Range("J113").Value = 1672.25
$ws.Range("K113").Value = 1204.5294
$ws.Range("L113").Value = 1672.25
$ws.Range("M113").Value = 965.4706000000001
$ws.Range("N113").Value = -6012.25

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 6523.48
$ws.Range("I126").Value = 9511.714
$ws.Range("J126").Value = 5361.3887
$ws.Range("K126").Value = 28535.142
$ws.Range("L126").Value = 16084.1661
$ws.Range("M126").Value = -26065.142
$ws.Range("N126").Value = -21024.1661

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 8627069
$ws.Range("I132").Value = 3725.139
$ws.Range("J132").Value = 22737996
$ws.Range("K132").Value = 11175.417
$ws.Range("L132").Value = 68213988
$ws.Range("M132").Value = -8645.417000000001
$ws.Range("N132").Value = -68219048
